$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.911.29'
$ws.Range("E2").Value = '  +0.25%  '

$ws.Range("D3").Value = '1.890.66'
$ws.Range("E3").Value = '  -0.50%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.8210'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.48%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3220'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.17%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.40'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07019'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.35%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08028'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7451'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.90%  '

$ws.Range("D13").Value = '1.892.23'
$ws.Range("E13").Value = '  -0.50%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.191'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.14%  '

$ws.Range("E15").Value = '  +0.60%  '

$ws.Range("D16").Value = '29.909.24'
$ws.Range("E16").Value = '  +0.20%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.99'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.57%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.875'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.42%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.61'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.50%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007744'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.12%  '

$ws.Range("D22").Value = '2.144.18'
$ws.Range("E22").Value = '  -0.62%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.908'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.51%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1559'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +20.62%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.81'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.174'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.80'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.078'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.38%  '

$ws.Range("E30").Value = '  -1.83%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.518'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.55%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.268'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05624'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.99%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.069'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.07%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.269'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.82%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7290'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.10%  '

$ws.Range("E37").Value = '  -0.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01912'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.17%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.777'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.13%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4414'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '71.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.47%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.949'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.95%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8440'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.79%  '

$ws.Range("E44").Value = '  +0.06%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.871'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.42%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.77'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.43%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.567'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.53%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.682'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.84%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '988.08'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.56%  '

$ws.Range("D50").Value = '2.041.28'
$ws.Range("E50").Value = '  -0.68%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.96'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.57%  '

